# Update the Handback status timestamps on the zh-cn and de-de sheets.
# These values are stored as plain text (not Excel date serials), so we
# force a text assignment to avoid Excel auto-converting them to dates.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 08:19:46"
$wsZhCn.Range("E3").Value = "2016-03-21 08:19:46"
$wsZhCn.Range("H2").Value = "2016-03-21 08:20:10"
$wsZhCn.Range("H3").Value = "2016-03-21 08:20:10"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 08:19:50"
$wsDeDe.Range("E3").Value = "2016-03-21 08:19:50"
$wsDeDe.Range("H2").Value = "2016-03-21 08:20:16"
$wsDeDe.Range("H3").Value = "2016-03-21 08:20:16"
